# dialogue.xlsx: rework the "press a key" keybind strings.
#
# menuSubText (row 6) used to read "press any key to continue / restart";
# it now prompts the player to press the "S" key to start the game, and
# the Korean column instead carries the (quote-prefixed, curly-quoted)
# "press 'R' to restart" copy that used to live only in gameOverSubtext.
# gameOverSubtext (row 11) keeps its meaning but the Korean/Japanese
# copies switch from straight quotes around R to curly/corner quotes.
#
# All text is written via chained [char] codepoints (prefixed with an
# empty string literal so PowerShell concatenates instead of summing the
# char codes) to guarantee the exact Unicode characters land in the
# workbook regardless of this console's own text encoding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 : menuSubText -----------------------------------------------
# English: "Press the "S" key to start"
$ws.Range("B6").Value = "" + [char]0x50 + [char]0x72 + [char]0x65 + [char]0x73 + [char]0x73 + [char]0x20 + [char]0x74 + [char]0x68 + [char]0x65 + [char]0x20 + [char]0x22 + [char]0x53 + [char]0x22 + [char]0x20 + [char]0x6B + [char]0x65 + [char]0x79 + [char]0x20 + [char]0x74 + [char]0x6F + [char]0x20 + [char]0x73 + [char]0x74 + [char]0x61 + [char]0x72 + [char]0x74

# Czech: "Stiskněte klávesu “S” pro spuštění"
$ws.Range("C6").Value = "" + [char]0x53 + [char]0x74 + [char]0x69 + [char]0x73 + [char]0x6B + [char]0x6E + [char]0x11B + [char]0x74 + [char]0x65 + [char]0x20 + [char]0x6B + [char]0x6C + [char]0xE1 + [char]0x76 + [char]0x65 + [char]0x73 + [char]0x75 + [char]0x20 + [char]0x201C + [char]0x53 + [char]0x201D + [char]0x20 + [char]0x70 + [char]0x72 + [char]0x6F + [char]0x20 + [char]0x73 + [char]0x70 + [char]0x75 + [char]0x161 + [char]0x74 + [char]0x11B + [char]0x6E + [char]0xED

# Korean: "'R' 키를 눌러 시작하세요" - leading straight apostrophe marks this as
# an explicit quote-prefixed literal (Excel strips it from the stored text
# and instead records a "quotePrefix" cell style), then the real text
# starts with a curly left single quote.
$ws.Range("D6").Value = "'" + [char]0x2018 + [char]0x52 + [char]0x2019 + [char]0x20 + [char]0xD0A4 + [char]0xB97C + [char]0x20 + [char]0xB20C + [char]0xB7EC + [char]0x20 + [char]0xC2DC + [char]0xC791 + [char]0xD558 + [char]0xC138 + [char]0xC694

# Japanese: "「S」キーを押してスタートしてください"
$ws.Range("E6").Value = "" + [char]0x300C + [char]0x53 + [char]0x300D + [char]0x30AD + [char]0x30FC + [char]0x3092 + [char]0x62BC + [char]0x3057 + [char]0x3066 + [char]0x30B9 + [char]0x30BF + [char]0x30FC + [char]0x30C8 + [char]0x3057 + [char]0x3066 + [char]0x304F + [char]0x3060 + [char]0x3055 + [char]0x3044

# --- Row 11 : gameOverSubtext -------------------------------------------
# Korean: "다시 시작하려면 'R' 을 누르세요" (curly single quotes around R)
$ws.Range("D11").Value = "" + [char]0xB2E4 + [char]0xC2DC + [char]0x20 + [char]0xC2DC + [char]0xC791 + [char]0xD558 + [char]0xB824 + [char]0xBA74 + [char]0x20 + [char]0x2018 + [char]0x52 + [char]0x2019 + [char]0x20 + [char]0xC744 + [char]0x20 + [char]0xB204 + [char]0xB974 + [char]0xC138 + [char]0xC694

# Japanese: "再起動するには 「R」 を押してください" (corner brackets around R)
$ws.Range("E11").Value = "" + [char]0x518D + [char]0x8D77 + [char]0x52D5 + [char]0x3059 + [char]0x308B + [char]0x306B + [char]0x306F + [char]0x20 + [char]0x300C + [char]0x52 + [char]0x300D + [char]0x20 + [char]0x3092 + [char]0x62BC + [char]0x3057 + [char]0x3066 + [char]0x304F + [char]0x3060 + [char]0x3055 + [char]0x3044

# --- Sheet view: selection moved from B5 to C6 ---------------------------
$ws.Range("C6").Select()

# --- Page setup: explicit portrait orientation ---------------------------
$ws.PageSetup.Orientation = 1
